# Auto-generated edit script
# Commit: Add data for 2024-07-25
# Applies updated 2024 (and a few prior-year) crime-count values
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$edits = @(
    @{ Sheet = 'Citywide Totals'; Cell = 'K2'; Value = 4455 }
    @{ Sheet = 'Citywide Totals'; Cell = 'J3'; Value = 8080 }
    @{ Sheet = 'Citywide Totals'; Cell = 'K3'; Value = 4570 }
    @{ Sheet = 'Citywide Totals'; Cell = 'E4'; Value = 2034 }
    @{ Sheet = 'Citywide Totals'; Cell = 'I4'; Value = 1799 }
    @{ Sheet = 'Citywide Totals'; Cell = 'J4'; Value = 1827 }
    @{ Sheet = 'Citywide Totals'; Cell = 'K4'; Value = 921 }
    @{ Sheet = 'Citywide Totals'; Cell = 'K5'; Value = 330 }
    @{ Sheet = 'Citywide Totals'; Cell = 'I6'; Value = 8963 }
    @{ Sheet = 'Citywide Totals'; Cell = 'K6'; Value = 5147 }
    @{ Sheet = 'Citywide Totals'; Cell = 'E7'; Value = 26039 }
    @{ Sheet = 'Citywide Totals'; Cell = 'I7'; Value = 26257 }
    @{ Sheet = 'Citywide Totals'; Cell = 'J7'; Value = 29298 }
    @{ Sheet = 'Citywide Totals'; Cell = 'K7'; Value = 15423 }
    @{ Sheet = 'Logan Square'; Cell = 'K3'; Value = 49 }
    @{ Sheet = 'Logan Square'; Cell = 'K7'; Value = 205 }
    @{ Sheet = 'Austin'; Cell = 'K2'; Value = 288 }
    @{ Sheet = 'Austin'; Cell = 'K3'; Value = 310 }
    @{ Sheet = 'Austin'; Cell = 'K6'; Value = 347 }
    @{ Sheet = 'Austin'; Cell = 'K7'; Value = 1031 }
    @{ Sheet = 'South Chicago'; Cell = 'K2'; Value = 114 }
    @{ Sheet = 'South Chicago'; Cell = 'K7'; Value = 329 }
    @{ Sheet = 'Garfield Park'; Cell = 'K2'; Value = 179 }
    @{ Sheet = 'Garfield Park'; Cell = 'K3'; Value = 240 }
    @{ Sheet = 'Garfield Park'; Cell = 'K4'; Value = 29 }
    @{ Sheet = 'Garfield Park'; Cell = 'K6'; Value = 188 }
    @{ Sheet = 'Garfield Park'; Cell = 'K7'; Value = 650 }
    @{ Sheet = 'West Pullman'; Cell = 'K3'; Value = 93 }
    @{ Sheet = 'West Pullman'; Cell = 'K7'; Value = 274 }
    @{ Sheet = 'Grand Crossing'; Cell = 'K3'; Value = 173 }
    @{ Sheet = 'Grand Crossing'; Cell = 'J4'; Value = 34 }
    @{ Sheet = 'Grand Crossing'; Cell = 'K4'; Value = 23 }
    @{ Sheet = 'Grand Crossing'; Cell = 'J7'; Value = 903 }
    @{ Sheet = 'Grand Crossing'; Cell = 'K7'; Value = 523 }
    @{ Sheet = 'New City'; Cell = 'K2'; Value = 106 }
    @{ Sheet = 'New City'; Cell = 'K6'; Value = 138 }
    @{ Sheet = 'New City'; Cell = 'K7'; Value = 350 }
    @{ Sheet = 'Woodlawn'; Cell = 'K3'; Value = 107 }
    @{ Sheet = 'Woodlawn'; Cell = 'K4'; Value = 15 }
    @{ Sheet = 'Woodlawn'; Cell = 'K6'; Value = 62 }
    @{ Sheet = 'Woodlawn'; Cell = 'K7'; Value = 259 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K4'; Value = 59 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K7'; Value = 454 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K8'; Value = 1031 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K11'; Value = 305 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K14'; Value = 87 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K16'; Value = 52 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K19'; Value = 466 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K20'; Value = 355 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K23'; Value = 157 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K27'; Value = 141 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K29'; Value = 818 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K31'; Value = 169 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K33'; Value = 650 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J37'; Value = 903 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K37'; Value = 523 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K42'; Value = 571 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K47'; Value = 97 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K49'; Value = 89 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K52'; Value = 407 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K53'; Value = 205 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K54'; Value = 288 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K55'; Value = 173 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K60'; Value = 100 }
    @{ Sheet = 'By Neighborhood'; Cell = 'E63'; Value = 369 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I63'; Value = 219 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J63'; Value = 112 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K63'; Value = 46 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K65'; Value = 350 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K66'; Value = 50 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K67'; Value = 594 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K73'; Value = 136 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K75'; Value = 54 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K76'; Value = 213 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K77'; Value = 110 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K78'; Value = 183 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K80'; Value = 52 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K83'; Value = 329 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K84'; Value = 111 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K85'; Value = 692 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K86'; Value = 103 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K88'; Value = 179 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K93'; Value = 58 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K94'; Value = 190 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K95'; Value = 274 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K96'; Value = 172 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K99'; Value = 259 }
    @{ Sheet = 'By Neighborhood'; Cell = 'E101'; Value = 26039 }
    @{ Sheet = 'By Neighborhood'; Cell = 'I101'; Value = 26257 }
    @{ Sheet = 'By Neighborhood'; Cell = 'J101'; Value = 29298 }
    @{ Sheet = 'By Neighborhood'; Cell = 'K101'; Value = 15423 }
    @{ Sheet = 'Gage Park'; Cell = 'K6'; Value = 59 }
    @{ Sheet = 'Gage Park'; Cell = 'K7'; Value = 169 }
    @{ Sheet = 'North Lawndale'; Cell = 'K3'; Value = 210 }
    @{ Sheet = 'North Lawndale'; Cell = 'K7'; Value = 594 }
    @{ Sheet = 'South Deering'; Cell = 'K2'; Value = 33 }
    @{ Sheet = 'South Deering'; Cell = 'K7'; Value = 111 }
    @{ Sheet = 'Lincoln Park'; Cell = 'K6'; Value = 51 }
    @{ Sheet = 'Lincoln Park'; Cell = 'K7'; Value = 89 }
    @{ Sheet = 'Loop'; Cell = 'K6'; Value = 147 }
    @{ Sheet = 'Loop'; Cell = 'K7'; Value = 288 }
    @{ Sheet = 'Englewood'; Cell = 'K2'; Value = 235 }
    @{ Sheet = 'Englewood'; Cell = 'K3'; Value = 290 }
    @{ Sheet = 'Englewood'; Cell = 'K5'; Value = 24 }
    @{ Sheet = 'Englewood'; Cell = 'K6'; Value = 228 }
    @{ Sheet = 'Englewood'; Cell = 'K7'; Value = 818 }
    @{ Sheet = 'Chatham'; Cell = 'K2'; Value = 145 }
    @{ Sheet = 'Chatham'; Cell = 'K7'; Value = 466 }
    @{ Sheet = 'River North'; Cell = 'K6'; Value = 118 }
    @{ Sheet = 'River North'; Cell = 'K7'; Value = 213 }
    @{ Sheet = 'Bridgeport'; Cell = 'K2'; Value = 31 }
    @{ Sheet = 'Bridgeport'; Cell = 'K7'; Value = 87 }
    @{ Sheet = 'Humboldt Park'; Cell = 'K3'; Value = 176 }
    @{ Sheet = 'Humboldt Park'; Cell = 'K6'; Value = 215 }
    @{ Sheet = 'Humboldt Park'; Cell = 'K7'; Value = 571 }
    @{ Sheet = 'Rogers Park'; Cell = 'K4'; Value = 18 }
    @{ Sheet = 'Rogers Park'; Cell = 'K6'; Value = 67 }
    @{ Sheet = 'Rogers Park'; Cell = 'K7'; Value = 183 }
    @{ Sheet = 'Lower West Side'; Cell = 'K3'; Value = 47 }
    @{ Sheet = 'Lower West Side'; Cell = 'K7'; Value = 173 }
    @{ Sheet = 'Douglas'; Cell = 'K5'; Value = 7 }
    @{ Sheet = 'Douglas'; Cell = 'K7'; Value = 157 }
    @{ Sheet = 'West Ridge'; Cell = 'K2'; Value = 55 }
    @{ Sheet = 'West Ridge'; Cell = 'K7'; Value = 172 }
    @{ Sheet = 'Chicago Lawn'; Cell = 'K2'; Value = 122 }
    @{ Sheet = 'Chicago Lawn'; Cell = 'K3'; Value = 112 }
    @{ Sheet = 'Chicago Lawn'; Cell = 'K4'; Value = 11 }
    @{ Sheet = 'Chicago Lawn'; Cell = 'K6'; Value = 104 }
    @{ Sheet = 'Chicago Lawn'; Cell = 'K7'; Value = 355 }
    @{ Sheet = 'West Lawn'; Cell = 'K4'; Value = 3 }
    @{ Sheet = 'West Lawn'; Cell = 'K7'; Value = 58 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'K2'; Value = 159 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'K3'; Value = 143 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'K6'; Value = 116 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'K7'; Value = 454 }
    @{ Sheet = 'West Loop'; Cell = 'K3'; Value = 34 }
    @{ Sheet = 'West Loop'; Cell = 'K7'; Value = 190 }
    @{ Sheet = 'Kenwood'; Cell = 'K3'; Value = 29 }
    @{ Sheet = 'Kenwood'; Cell = 'K7'; Value = 97 }
    @{ Sheet = 'North Center'; Cell = 'K6'; Value = 26 }
    @{ Sheet = 'North Center'; Cell = 'K7'; Value = 50 }
    @{ Sheet = 'Belmont Cragin'; Cell = 'K3'; Value = 77 }
    @{ Sheet = 'Belmont Cragin'; Cell = 'K4'; Value = 17 }
    @{ Sheet = 'Belmont Cragin'; Cell = 'K6'; Value = 109 }
    @{ Sheet = 'Belmont Cragin'; Cell = 'K7'; Value = 305 }
    @{ Sheet = 'Portage Park'; Cell = 'K3'; Value = 33 }
    @{ Sheet = 'Portage Park'; Cell = 'K7'; Value = 136 }
    @{ Sheet = 'United Center'; Cell = 'K6'; Value = 74 }
    @{ Sheet = 'United Center'; Cell = 'K7'; Value = 179 }
    @{ Sheet = 'Edgewater'; Cell = 'K3'; Value = 35 }
    @{ Sheet = 'Edgewater'; Cell = 'K7'; Value = 141 }
    @{ Sheet = 'Streeterville'; Cell = 'K4'; Value = 40 }
    @{ Sheet = 'Streeterville'; Cell = 'K7'; Value = 103 }
    @{ Sheet = 'Pullman'; Cell = 'K2'; Value = 20 }
    @{ Sheet = 'Pullman'; Cell = 'K7'; Value = 54 }
    @{ Sheet = 'Morgan Park'; Cell = 'K6'; Value = 26 }
    @{ Sheet = 'Morgan Park'; Cell = 'K7'; Value = 100 }
    @{ Sheet = 'South Shore'; Cell = 'K3'; Value = 231 }
    @{ Sheet = 'South Shore'; Cell = 'K4'; Value = 40 }
    @{ Sheet = 'South Shore'; Cell = 'K7'; Value = 692 }
    @{ Sheet = 'Riverdale'; Cell = 'K2'; Value = 47 }
    @{ Sheet = 'Riverdale'; Cell = 'K7'; Value = 110 }
    @{ Sheet = 'Rush & Division'; Cell = 'K3'; Value = 11 }
    @{ Sheet = 'Rush & Division'; Cell = 'K7'; Value = 52 }
    @{ Sheet = 'Little Village'; Cell = 'K2'; Value = 106 }
    @{ Sheet = 'Little Village'; Cell = 'K7'; Value = 407 }
    @{ Sheet = 'Archer Heights'; Cell = 'K2'; Value = 20 }
    @{ Sheet = 'Archer Heights'; Cell = 'K7'; Value = 59 }
    @{ Sheet = 'Bucktown'; Cell = 'K6'; Value = 30 }
    @{ Sheet = 'Bucktown'; Cell = 'K7'; Value = 52 }
)

foreach ($edit in $edits) {
    $ws = $wb.Worksheets.Item($edit.Sheet)
    $ws.Range($edit.Cell).Value = $edit.Value
}

Write-Host "Applied $($edits.Count) cell updates."
